# hydro_cascades.xlsx — fix Excel mistake with parameter kappaYieldH
#
# Column L ("kappaYieldH" style power/yield calc) was wired to the wrong
# source columns (I*J instead of J*K). Re-point every formula in L3:L42
# at J/K of the same row, then restore the active cell/selection that was
# left on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 42; $r++) {
    $ws.Cells.Item($r, 12).Formula = "=J$r*K$r*9.8/1000"
}

# Restore selection left on the sheet after the edit (was M4, now L44)
$ws.Range("L44").Select()
